$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing row 6 (which already carries the
# correct style indices) down into the still-empty row 5 so the newly
# typed data picks up matching cell styles.
$ws.Range("A6:J6").Copy() | Out-Null
$ws.Range("A5:J5").PasteSpecial(-4122) | Out-Null

# Row 5: "Crumpet GEF" facility
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# Row 6: "Scone GEF" facility
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# Leave the sheet with the same selection the author ended up with.
$ws.Range("A5:J6").Select() | Out-Null
